$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.179.84'
$ws.Range('E2').Value = '  -0.41%  '
$ws.Range('D3').Value = '2.078.10'
$ws.Range('E3').Value = '  -0.75%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '''253.97'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.97%  '
$ws.Range('E6').Value = '  +1.99%  '
$ws.Range('D7').Value = '''59.80'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +10.46%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '''0.394'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +5.03%  '
$ws.Range('D10').Value = '''61.51'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.44%  '
$ws.Range('E11').Value = '  +8.05%  '
$ws.Range('E12').Value = '  +2.62%  '
$ws.Range('D13').Value = '''16.32'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +7.67%  '
$ws.Range('D14').Value = '2.379.85'
$ws.Range('E14').Value = '  -0.82%  '
$ws.Range('E15').Value = '  -2.08%  '
$ws.Range('D16').Value = '''5.60'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +8.21%  '
$ws.Range('D17').Value = '2.080.34'
$ws.Range('E17').Value = '  -0.88%  '
$ws.Range('D18').Value = '37.286.96'
$ws.Range('E18').Value = '  +0.11%  '
$ws.Range('D19').Value = '''16.18'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +10.97%  '
$ws.Range('D20').Value = '''74.89'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.91%  '
$ws.Range('D21').Value = '0.0₃0930'
$ws.Range('E21').Value = '  +9.71%  '
$ws.Range('D22').Value = '''5.50'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.78%  '
$ws.Range('D23').Value = '''239.66'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.53%  '
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('E25').Value = '  -2.22%  '
$ws.Range('D26').Value = '''2.30'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +14.75%  '
$ws.Range('D27').Value = '''170.43'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E28').Value = '  +1.70%  '
$ws.Range('D29').Value = '''20.43'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.25%  '
$ws.Range('E30').Value = '  +3.08%  '
$ws.Range('E31').Value = '  +6.51%  '
$ws.Range('D33').Value = '''0.0639'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.51%  '
$ws.Range('E34').Value = '  +8.95%  '
$ws.Range('D35').Value = '''0.0914'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('E36').Value = '  -0.15%  '
$ws.Range('E37').Value = '  +2.65%  '
$ws.Range('D38').Value = '''0.119'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +27.66%  '
$ws.Range('E39').Value = '  -4.35%  '
$ws.Range('E40').Value = '  +2.20%  '
$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D41').Value = '''17.96'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.89%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').Value = '''0.0228'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.55%  '
$ws.Range('E43').Value = '  +0.06%  '
$ws.Range('D44').Value = '''99.31'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.53%  '
$ws.Range('D45').Value = '''4.40'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.75%  '
$ws.Range('D46').Value = '''2.86'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.69%  '
$ws.Range('D47').Value = '''4.60'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +14.04%  '
$ws.Range('E48').Value = '  +7.66%  '
$ws.Range('D49').Value = '1.311.27'
$ws.Range('E49').Value = '  -0.89%  '
$ws.Range('E50').Value = '  +0.22%  '
$ws.Range('D51').Value = '''6.97'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.45%  '
